# Add a "Status" / "FAIL" column to the credentials sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell C1: "Status", styled like the other header cells (bold font
# + yellow fill) but without the border the existing headers have.
$ws.Range("C1").Value = "Status"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").Interior.Color = 65535

# Data cell C2: "FAIL", left with default formatting.
$ws.Range("C2").Value = "FAIL"

# Leave the selection on H7, matching the saved workbook view.
$ws.Range("H7").Select() | Out-Null
